$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading. All the edits below
# only apply to bullets that live after that heading, because a couple of
# the bullet strings also appear verbatim earlier in the document (under
# "Partner - Siege Analytics"), and those earlier copies must stay as-is.
$headingRange = $d.Content
$headingRange.Find.Execute("KEY ACHIEVEMENTS AND IMPACT", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$sectionStart = $headingRange.End

function Get-SectionRange() {
    return $d.Range($sectionStart, $d.Content.End)
}

function Replace-InSection($oldText, $newText) {
    $r = Get-SectionRange
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

function Remove-BulletInSection($text) {
    $r = Get-SectionRange
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        # Extend the found range by one character to swallow the trailing
        # paragraph mark too, so the whole bullet paragraph disappears.
        $d.Range($r.Start, $r.End + 1).Delete()
    }
}

# 1) "Discovered systematic race coding errors ..." -> predictive excellence bullet
Replace-InSection `
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%" `
    "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"

# 2) "Achieved 87% prediction accuracy ..." -> reduced polling margins bullet
Replace-InSection `
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%" `
    "• Reduced polling margins from ±4.2% to ±2.1%"

# 3) "Built redistricting platform ..." -> executive authority bullet
Replace-InSection `
    "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations" `
    "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"

# 4) "Developed longitudinal data analysis methods ..." -> methodological advancement bullet
Replace-InSection `
    "• Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality" `
    "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"

# 5) Drop the trailing "Provided expert testimony ..." bullet entirely
Remove-BulletInSection "• Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy"

# 6) Drop the trailing "Demystified FEC compliance ..." bullet entirely
Remove-BulletInSection "• Demystified FEC compliance through real-time processing systems enabling transparent campaign finance monitoring"
